# Apply cryptos list update (auto-generated from diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "219.41").
# Force those cells to remain plain text so the exact string is preserved,
# matching the original inlineStr cell content, then restore the default style.
$priceRows = @(2, 3, 5, 6, 7, 8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 22, 23, 25, 26, 27, 33, 36, 41, 42, 45, 47, 48, 49, 50, 51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "28.286.25"
$ws.Cells.Item(2, 5).Value = "  +4.20%  "
$ws.Cells.Item(3, 4).Value = "1.727.45"
$ws.Cells.Item(3, 5).Value = "  +2.82%  "
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$ws.Cells.Item(5, 4).Value = "219.41"
$ws.Cells.Item(5, 5).Value = "  +1.92%  "
$ws.Cells.Item(6, 4).Value = "0.522"
$ws.Cells.Item(6, 5).Value = "  +0.68%  "
$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(8, 4).Value = "24.40"
$ws.Cells.Item(8, 5).Value = "  +14.31%  "
$ws.Cells.Item(9, 4).Value = "0.265"
$ws.Cells.Item(9, 5).Value = "  +3.55%  "
$ws.Cells.Item(10, 5).Value = "  +1.88%  "
$ws.Cells.Item(11, 4).Value = "0.0899"
$ws.Cells.Item(11, 5).Value = "  +1.35%  "
$ws.Cells.Item(12, 4).Value = "1.970.20"
$ws.Cells.Item(12, 5).Value = "  +2.75%  "
$ws.Cells.Item(13, 4).Value = "1.728.75"
$ws.Cells.Item(13, 5).Value = "  +2.92%  "
$ws.Cells.Item(14, 4).Value = "4.28"
$ws.Cells.Item(14, 5).Value = "  +3.24%  "
$ws.Cells.Item(15, 4).Value = "0.563"
$ws.Cells.Item(15, 5).Value = "  +4.83%  "
$ws.Cells.Item(16, 4).Value = "67.71"
$ws.Cells.Item(16, 5).Value = "  +2.22%  "
$ws.Cells.Item(17, 4).Value = "28.230.69"
$ws.Cells.Item(17, 5).Value = "  +4.02%  "
$ws.Cells.Item(18, 4).Value = "244.43"
$ws.Cells.Item(18, 5).Value = "  +2.15%  "
$ws.Cells.Item(19, 4).Value = "8.03"
$ws.Cells.Item(19, 5).Value = "  -0.38%  "
$ws.Cells.Item(20, 5).Value = "  +1.87%  "
$ws.Cells.Item(21, 5).Value = "  -0.13%  "
$ws.Cells.Item(22, 4).Value = "4.65"
$ws.Cells.Item(22, 5).Value = "  +2.87%  "
$ws.Cells.Item(23, 4).Value = "9.70"
$ws.Cells.Item(23, 5).Value = "  +2.60%  "
$ws.Cells.Item(24, 5).Value = "  -0.35%  "
$ws.Cells.Item(25, 4).Value = "149.28"
$ws.Cells.Item(25, 5).Value = "  +1.50%  "
$ws.Cells.Item(26, 4).Value = "7.53"
$ws.Cells.Item(26, 5).Value = "  +3.93%  "
$ws.Cells.Item(27, 4).Value = "16.75"
$ws.Cells.Item(27, 5).Value = "  +2.33%  "
$ws.Cells.Item(28, 5).Value = "  +0.94%  "
$ws.Cells.Item(29, 5).Value = "  -0.30%  "
$ws.Cells.Item(30, 5).Value = "  +2.73%  "
$ws.Cells.Item(31, 5).Value = "  +2.21%  "
$ws.Cells.Item(32, 5).Value = "  +2.20%  "
$ws.Cells.Item(33, 4).Value = "1.504.53"
$ws.Cells.Item(33, 5).Value = "  -3.83%  "
$ws.Cells.Item(34, 5).Value = "  +2.07%  "
$ws.Cells.Item(35, 5).Value = "  -1.28%  "
$ws.Cells.Item(36, 4).Value = "0.966"
$ws.Cells.Item(36, 5).Value = "  +3.70%  "
$ws.Cells.Item(37, 5).Value = "  +1.63%  "
$ws.Cells.Item(38, 5).Value = "  +0.46%  "
$ws.Cells.Item(39, 5).Value = "  +0.99%  "
$ws.Cells.Item(40, 5).Value = "  +1.19%  "
$ws.Cells.Item(41, 4).Value = "71.07"
$ws.Cells.Item(41, 5).Value = "  +2.76%  "
$ws.Cells.Item(42, 4).Value = "5.78"
$ws.Cells.Item(42, 5).Value = "  +3.94%  "
$ws.Cells.Item(43, 5).Value = "  -0.19%  "
$ws.Cells.Item(44, 5).Value = "  +2.07%  "
$ws.Cells.Item(45, 4).Value = "1.875.99"
$ws.Cells.Item(45, 5).Value = "  +2.65%  "
$ws.Cells.Item(46, 5).Value = "  +2.77%  "
$ws.Cells.Item(47, 4).Value = "1.77"
$ws.Cells.Item(47, 5).Value = "  +11.95%  "
$ws.Cells.Item(48, 4).Value = "91.02"
$ws.Cells.Item(48, 5).Value = "  +0.41%  "
$ws.Cells.Item(49, 4).Value = "0.0₆0114"
$ws.Cells.Item(49, 5).Value = "  +6.50%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "8.26"
$ws.Cells.Item(50, 5).Value = "  +2.12%  "
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).Value = "0.105"
$ws.Cells.Item(51, 5).Value = "  +1.35%  "

foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
